# Fruta / hortaliza, semanal
#
# A new weekly price record needs to be inserted for "Macroferia Regional de
# Talca - Kiwi" ahead of the existing row 180, pushing the former rows
# 180-202 down to 181-203 (dimension grows from A1:T202 to A1:T203).
# The new row carries the same fixed attributes (market, region, product
# taxonomy, province) as its neighbours, with its own date / grade / price
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 180; Excel shifts rows 180:202
# down to 181:203 and copies formatting (incl. the date style on column D)
# from the row above, matching the rest of the table.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new weekly record.
$ws.Cells.Item(180, 1).Value = 5
$ws.Cells.Item(180, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(180, 3).Value = "Maule"
$ws.Cells.Item(180, 4).Value = 44504
$ws.Cells.Item(180, 5).Value = 7
$ws.Cells.Item(180, 6).Value = "Fruta"
$ws.Cells.Item(180, 7).Value = 100101
$ws.Cells.Item(180, 8).Value = "Berries"
$ws.Cells.Item(180, 9).Value = 100101007
$ws.Cells.Item(180, 10).Value = "Kiwi"
$ws.Cells.Item(180, 11).Value = "Hayward"
$ws.Cells.Item(180, 12).Value = "Primera"
$ws.Cells.Item(180, 13).Value = 80
$ws.Cells.Item(180, 14).Value = 16000
$ws.Cells.Item(180, 15).Value = 16000
$ws.Cells.Item(180, 16).Value = 16000
$ws.Cells.Item(180, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(180, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(180, 19).Value = 889
$ws.Cells.Item(180, 20).Value = 18
